$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 6562.625
$ws.Range("I32").Value = 1500.5
$ws.Range("J32").Value = 8250
$ws.Range("K32").Value = 1500.5
$ws.Range("L32").Value = 8250
$ws.Range("M32").Value = -1174.5
$ws.Range("N32").Value = -8902
$ws.Range("H33").Value = 369
$ws.Range("J33").Value = 575
$ws.Range("L33").Value = 575
$ws.Range("N33").Value = -1033
$ws.Range("H103").Value = 990
$ws.Range("J103").Value = 990
$ws.Range("L103").Value = 2970
$ws.Range("N103").Value = -4142
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()
$ws.Range("H130").Value = 96385
$ws.Range("J130").Value = 96385
$ws.Range("L130").Value = 96385
$ws.Range("N130").Value = -106425
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
$ws.Range("H137").Value = 3126.1333
$ws.Range("I137").Value = 2645.6924
$ws.Range("K137").Value = 7937.0772
$ws.Range("M137").Value = -5387.0772

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()
$ws.Range("H122").Value = 2194.625
$ws.Range("I122").Value = 2258.3333
$ws.Range("J122").Value = 2112.7144
$ws.Range("K122").Value = 6774.999899999999
$ws.Range("L122").Value = 6338.1432
$ws.Range("M122").Value = -4324.999899999999
$ws.Range("N122").Value = -11238.1432
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()
$ws.Range("H133").Value = 99994
$ws.Range("J133").Value = 99994
$ws.Range("L133").Value = 99994
$ws.Range("N133").Value = -105054
$ws.Range("H134").Value = 96000
$ws.Range("J134").Value = 96000
$ws.Range("L134").Value = 96000
$ws.Range("N134").Value = -106140
$ws.Range("H135").Value = 89000
$ws.Range("J135").Value = 89000
$ws.Range("L135").Value = 89000
$ws.Range("N135").Value = -99140

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()
$ws.Range("H115").Value = 0
$ws.Range("I115").Value = 0
$ws.Range("K115").Value = 0
$ws.Range("M115").ClearContents()
$ws.Range("H129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").ClearContents()

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6978.0835
$ws.Range("I31").Value = 2991
$ws.Range("K31").Value = 2991
$ws.Range("M31").Value = -2696
$ws.Range("H34").Value = 6978.0835
$ws.Range("I34").Value = 2991
$ws.Range("K34").Value = 2991
$ws.Range("M34").Value = -2789
$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("M58").ClearContents()
$ws.Range("N58").ClearContents()
$ws.Range("H60").Value = 17016.666
$ws.Range("I60").Value = 22000
$ws.Range("J60").Value = 16563.637
$ws.Range("K60").Value = 22000
$ws.Range("L60").Value = 16563.637
$ws.Range("M60").Value = -21489
$ws.Range("N60").Value = -17585.637
$ws.Range("H98").Value = 0
$ws.Range("I98").Value = 0
$ws.Range("K98").Value = 0
$ws.Range("M98").ClearContents()
$ws.Range("H110").Value = 99995
$ws.Range("J110").Value = 99995
$ws.Range("L110").Value = 99995
$ws.Range("N110").Value = -108175
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()
$ws.Range("H132").Value = 1728.125
$ws.Range("I132").Value = 1728.125
$ws.Range("K132").Value = 5184.375
$ws.Range("M132").Value = -2654.375
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("M136").ClearContents()
$ws.Range("N136").ClearContents()

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 749.75
$ws.Range("J107").Value = 833
$ws.Range("L107").Value = 2499
$ws.Range("N107").Value = -6339

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H103").Value = 49999
$ws.Range("J103").Value = 49999
$ws.Range("L103").Value = 49999
$ws.Range("N103").Value = -52343
$ws.Range("H110").Value = 99980
$ws.Range("J110").Value = 99980
$ws.Range("L110").Value = 99980
$ws.Range("N110").Value = -108160
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()
$ws.Range("H132").Value = 3797.5
$ws.Range("I132").Value = 3797.5
$ws.Range("K132").Value = 11392.5
$ws.Range("M132").Value = -8862.5

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3748.5
$ws.Range("I7").Value = 1664.6666
$ws.Range("K7").Value = 1664.6666
$ws.Range("M7").Value = -1552.6666
$ws.Range("H123").Value = 0
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("M123").ClearContents()
$ws.Range("N123").ClearContents()
$ws.Range("H126").Value = 3748.5
$ws.Range("I126").Value = 1664.6666
$ws.Range("K126").Value = 4993.9998
$ws.Range("M126").Value = -2523.9998
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()
$ws.Range("H140").Value = 72000
$ws.Range("J140").Value = 72000
$ws.Range("L140").Value = 72000
$ws.Range("N140").Value = -82360

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H104").Value = 21999.5
$ws.Range("J104").Value = 21999.5
$ws.Range("L104").Value = 21999.5
$ws.Range("N104").Value = -28987.5
$ws.Range("H113").Value = 1249.1666
$ws.Range("I113").Value = 480
$ws.Range("J113").Value = 2018.3334
$ws.Range("K113").Value = 1440
$ws.Range("L113").Value = 6055.0002
$ws.Range("M113").Value = 730
$ws.Range("N113").Value = -10395.0002
$ws.Range("H126").Value = 1165
$ws.Range("I126").Value = 997.5
$ws.Range("K126").Value = 2992.5
$ws.Range("M126").Value = -522.5
$ws.Range("H136").Value = 2193
$ws.Range("I136").Value = 2193
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 6579
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -4029
$ws.Range("N136").ClearContents()
$ws.Range("H137").Value = 39000
$ws.Range("J137").Value = 39000
$ws.Range("L137").Value = 39000
$ws.Range("N137").Value = -49200
Write-Host "All edits applied"
